# Refresh the "cryptos" price/volume table (GitHub Actions style scrape update).
# Price (D) and Volume(1h) (E) columns are plain text cells (t="inlineStr" in the
# original sheet) -- not real numbers -- so we assign string literals throughout.
# A handful of Price values (rows 32, 37, 38, 47, 51) have a trailing zero
# (e.g. "0.0890", "4.60", "0.100"); Excel's COM layer auto-coerces a bare numeric
# -looking string into a real number on assignment, which would silently drop
# that trailing zero. For those specific cells we prefix the literal with a
# leading apostrophe (the normal Excel "force text" convention) so the value is
# stored verbatim as text instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.754.45'
$ws.Range('E2').Value = '  -1.56%  '

$ws.Range('D3').Value = '2.350.19'
$ws.Range('E3').Value = '  -2.23%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').Value = '321.62'
$ws.Range('E5').Value = '  -1.06%  '

$ws.Range('D6').Value = '105.61'
$ws.Range('E6').Value = '  +0.40%  '

$ws.Range('E7').Value = '  -3.07%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').Value = '0.616'
$ws.Range('E9').Value = '  -6.68%  '

$ws.Range('D10').Value = '41.27'
$ws.Range('E10').Value = '  -2.52%  '

$ws.Range('E11').Value = '  -2.41%  '

$ws.Range('D12').Value = '8.42'
$ws.Range('E12').Value = '  -2.38%  '

$ws.Range('E13').Value = '  -2.31%  '

$ws.Range('E14').Value = '  -0.20%  '

$ws.Range('D15').Value = '15.98'
$ws.Range('E15').Value = '  -7.49%  '

$ws.Range('D16').Value = '2.706.62'
$ws.Range('E16').Value = '  -2.29%  '

$ws.Range('D17').Value = '2.341.59'
$ws.Range('E17').Value = '  -2.68%  '

$ws.Range('D18').Value = '42.760.81'

$ws.Range('D19').Value = '7.76'
$ws.Range('E19').Value = '  +4.31%  '

$ws.Range('E20').Value = '  -3.13%  '

$ws.Range('D21').Value = '77.25'
$ws.Range('E21').Value = '  +1.56%  '

$ws.Range('E22').Value = '  +3.80%  '

$ws.Range('D23').Value = '260.07'
$ws.Range('E23').Value = '  -4.67%  '

$ws.Range('D24').Value = '2.32'
$ws.Range('E24').Value = '  -4.67%  '

$ws.Range('D25').Value = '9.53'
$ws.Range('E25').Value = '  -2.04%  '

$ws.Range('E26').Value = '  +0.04%  '

$ws.Range('D27').Value = '11.38'
$ws.Range('E27').Value = '  -4.71%  '

$ws.Range('D28').Value = '23.07'
$ws.Range('E28').Value = '  +0.57%  '

$ws.Range('E29').Value = '  +0.25%  '

$ws.Range('D30').Value = '174.71'
$ws.Range('E30').Value = '  -1.95%  '

$ws.Range('D31').Value = '36.41'
$ws.Range('E31').Value = '  -4.33%  '

$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '''0.0890'
$ws.Range('E32').Value = '  -4.95%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '6.14'
$ws.Range('E33').Value = '  +3.26%  '

$ws.Range('E34').Value = '  -7.63%  '

$ws.Range('E35').Value = '  +8.76%  '

$ws.Range('E36').Value = '  -3.12%  '

$ws.Range('D37').Value = '''4.60'
$ws.Range('E37').Value = '  -5.79%  '

$ws.Range('D38').Value = '''0.0360'
$ws.Range('E38').Value = '  -2.88%  '

$ws.Range('D39').Value = '3.79'
$ws.Range('E39').Value = '  -7.69%  '

$ws.Range('E40').Value = '  -5.47%  '

$ws.Range('D41').Value = '71.81'
$ws.Range('E41').Value = '  +2.63%  '

$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = '0.235'
$ws.Range('E42').Value = '  +0.10%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '1.47'
$ws.Range('E43').Value = '  -8.48%  '

$ws.Range('E44').Value = '  -0.11%  '

$ws.Range('D45').Value = '115.73'
$ws.Range('E45').Value = '  -8.97%  '

$ws.Range('D46').Value = '88.78'
$ws.Range('E46').Value = '  -0.37%  '

$ws.Range('D47').Value = '''11.90'
$ws.Range('E47').Value = '  -6.77%  '

$ws.Range('E48').Value = '  -3.50%  '

$ws.Range('D49').Value = '9.15'
$ws.Range('E49').Value = '  -5.99%  '

$ws.Range('D50').Value = '73.65'
$ws.Range('E50').Value = '  +1.04%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.100'
$ws.Range('E51').Value = '  -1.46%  '
